# EnterpriseJavaTimeLog.xlsx — add Week 8 log entries.
#
# The existing sheet had a block of blank placeholder rows (49-58) below the
# last filled-in week (week 7, ending row 48) that were reserved for future
# weeks, followed by a "TODO / loose ends" section starting at row 55.
#
# This edit:
#   1) Inserts 4 blank rows before row 55 to make room (pushing the TODO
#      section, which used to start at row 55, down to row 59 and beyond).
#   2) Fills in the first three of the now-available blank rows (49-51) with
#      the Week 8 log entries: overview/reading, activities, and the final
#      notes for the week (including an indie-project note).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new entries by pushing the TODO/loose-ends block (and
# everything after it) down by 4 rows.
$ws.Rows("55:58").Insert()

# Week 8 overview video and reading -- 2019-03-25, 1.5 hrs
$ws.Range("A49").Value = 43549
$ws.Range("B49").Value = 1.5
$ws.Range("D49").Value = "Week 8 overview video and reading"

# Week 8 activities -- 2019-03-26, 5.5 hrs
$ws.Range("A50").Value = 43550
$ws.Range("B50").Value = 5.5
$ws.Range("D50").Value = "'Week 8 Activities (including time spent sorting out how to handle file permissions in jdk lib)"

# Week 8 activities (last things) + indie project note -- 2019-03-28, 3 hrs
$ws.Range("A51").Value = 43552
$ws.Range("B51").Value = 3
$ws.Range("D51").Value = "'Week 8 Activities (last things)`nIndie Project: tried to investigate template security issue… it's the package-lock.json.  Wonder if I could use the template without the javascript and json since it's mostly for layout??  Decided against."
$ws.Rows(51).RowHeight = 45

$ws.Range("D52").Select()
